$wb = $excel.ActiveWorkbook

# "About" sheet: clear the stray date value in C1 (leftover from a previous save)
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Clear()

# "Set Values Here" sheet: update the weight values in rows 9, 15 and 16
$wsSet = $wb.Worksheets.Item("Set Values Here")

# Row 9 (Deficit Spending weights for fuel tax revenue)
$wsSet.Range("C9").Value = 5
$wsSet.Range("D9").Value = 0
$wsSet.Range("F9").Value = 5

# Row 15 (Regular Spending weight for national debt interest)
$wsSet.Range("B15").Value = 5

# Row 16 (Regular Spending weight for remainder)
$wsSet.Range("B16").Value = 5
